$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 17

$ws.Cells.Item($row, 1).Value = "2025-08-20 07:56:50"
$ws.Cells.Item($row, 2).Value = "create-team"
$ws.Cells.Item($row, 3).Value = "new-organization97"
$ws.Cells.Item($row, 4).Value = "newteam"
$ws.Cells.Item($row, 5).Value = "demo"
$cI = $ws.Cells.Item($row, 9)
$cI.Value = "'False"
$cI.Style = "Normal"
